$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 18:50"

# Country data updates scraped at 18:50 (totals, new cases, active, recovered, critical, new deaths, deaths)

# Row 4: Estados Unidos
$ws.Range("B4").Value = 145542
$ws.Range("C4").Value = 2051
$ws.Range("D4").Value = 4579
$ws.Range("E4").Value = 138347
$ws.Range("G4").Value = 33
$ws.Range("H4").Value = 2616

# Row 12: Suiza
$ws.Range("B12").Value = 15760
$ws.Range("C12").Value = 931
$ws.Range("E12").Value = 13589
$ws.Range("G12").Value = 48
$ws.Range("H12").Value = 348

# Row 20: Noruega
$ws.Range("B20").Value = 4445
$ws.Range("C20").Value = 161
$ws.Range("E20").Value = 4401

# Row 21: Brasil
$ws.Range("B21").Value = 4371
$ws.Range("C21").Value = 115
$ws.Range("E21").Value = 4110

# Row 25: Chequia
$ws.Range("B25").Value = 2942
$ws.Range("C25").Value = 125
$ws.Range("E25").Value = 2914

# Row 32: Rumania
$ws.Range("D32").Value = 209
$ws.Range("E32").Value = 1696
$ws.Range("F32").Value = 33
$ws.Range("G32").Value = 4
$ws.Range("H32").Value = 47

# Row 58: Catar
$ws.Range("A58").Value = "Catar"
$ws.Range("B58").Value = 693
$ws.Range("C58").Value = 59
$ws.Range("D58").Value = 51
$ws.Range("E58").Value = 641
$ws.Range("F58").Value = 6
$ws.Range("H58").Value = 1

# Row 59: Hong Kong
$ws.Range("A59").Value = "Hong Kong"
$ws.Range("B59").Value = 642
$ws.Range("D59").Value = 118
$ws.Range("E59").Value = 520
$ws.Range("F59").Value = 5
$ws.Range("H59").Value = 4

# Row 63: Argelia
$ws.Range("D63").Value = 37
$ws.Range("E63").Value = 512

# Row 94: Afganistan
$ws.Range("A94").Value = "Afganistan"
$ws.Range("C94").Value = 50
$ws.Range("D94").Value = 2
$ws.Range("E94").Value = 164
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0

# Row 95: Cuba
$ws.Range("A95").Value = "Cuba"
$ws.Range("B95").Value = 170
$ws.Range("C95").Value = 31
$ws.Range("D95").Value = 4
$ws.Range("E95").Value = 162
$ws.Range("F95").Value = 2
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 4

# Row 96: Islas Feroe
$ws.Range("A96").Value = "Islas Feroe"
$ws.Range("B96").Value = 168
$ws.Range("C96").Value = 9
$ws.Range("D96").Value = 70
$ws.Range("E96").Value = 98
$ws.Range("F96").Value = 1
$ws.Range("H96").Value = 0

# Row 97: Costa de Marfil
$ws.Range("A97").Value = "Costa de Marfil"
$ws.Range("B97").Value = 165
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 4
$ws.Range("E97").Value = 160
$ws.Range("H97").Value = 1

# Row 98: Senegal
$ws.Range("A98").Value = "Senegal"
$ws.Range("B98").Value = 162
$ws.Range("C98").Value = 20
$ws.Range("D98").Value = 27
$ws.Range("E98").Value = 135
$ws.Range("F98").Value = 0
$ws.Range("H98").Value = 1

# Row 99: Malta
$ws.Range("A99").Value = "Malta"
$ws.Range("B99").Value = 156
$ws.Range("C99").Value = 5
$ws.Range("E99").Value = 154
$ws.Range("F99").Value = 4
$ws.Range("H99").Value = 0

# Row 100: Ghana
$ws.Range("A100").Value = "Ghana"
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 2
$ws.Range("E100").Value = 145
$ws.Range("F100").Value = 1
$ws.Range("H100").Value = 5

# Row 101: Bielorrusia
$ws.Range("A101").Value = "Bielorrusia"
$ws.Range("B101").Value = 152
$ws.Range("C101").Value = 58
$ws.Range("D101").Value = 32
$ws.Range("E101").Value = 120
$ws.Range("F101").Value = 2
$ws.Range("H101").Value = 0

# Row 102: Uzbekistan
$ws.Range("A102").Value = "Uzbekistan"
$ws.Range("B102").Value = 149
$ws.Range("C102").Value = 5
$ws.Range("D102").Value = 7
$ws.Range("E102").Value = 140
$ws.Range("F102").Value = 8
$ws.Range("H102").Value = 2

# Row 110: Nigeria
$ws.Range("E110").Value = 106
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 2

# Row 128: Monaco
$ws.Range("A128").Value = "Monaco"
$ws.Range("C128").Value = 3
$ws.Range("D128").Value = 1
$ws.Range("E128").Value = 47
$ws.Range("F128").Value = 0
$ws.Range("H128").Value = 1

# Row 129: Banglades
$ws.Range("A129").Value = "Banglades"
$ws.Range("B129").Value = 49
$ws.Range("C129").Value = 1
$ws.Range("D129").Value = 19
$ws.Range("E129").Value = 25
$ws.Range("F129").Value = 1
$ws.Range("H129").Value = 5
